# Generate Report for Handback
#
# Two files (307fbd6b-89cd-4f2d-882b-87102e101d0d and
# a69741f6-eb82-413a-a776-326854e80aea) move from "Ready for handoff" to
# "Handed back: in sync with en-US" for both the zh-cn and de-de locales.
# As part of the handback, the "Latest Target File" / "Latest Handback File"
# / "Latest Handback DateTime" columns on the per-locale sheets get filled
# in (they were previously blank / zero-date placeholders).

$wb = $excel.ActiveWorkbook

$READY_FOR_HANDOFF = "Ready for handoff"
$HANDED_BACK = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: update status columns (zh-cn / de-de) for rows 4 and 5
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E4").Value = $HANDED_BACK
$wsOverview.Range("F4").Value = $HANDED_BACK
$wsOverview.Range("E5").Value = $HANDED_BACK
$wsOverview.Range("F5").Value = $HANDED_BACK

# ---------------------------------------------------------------------
# zh-cn sheet: update status + fill in handback file / datetime columns
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Row 4 - 307fbd6b-89cd-4f2d-882b-87102e101d0d
$wsZhCn.Range("C4").Value = $HANDED_BACK
$wsZhCn.Range("I4").Value = "307fbd6b-89cd-4f2d-882b-87102e101d0d.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7c748e17e7c44fbaa2203baa9dbeb7529a5ddb67/e2e/307fbd6b-89cd-4f2d-882b-87102e101d0d.md") | Out-Null
$wsZhCn.Range("I4").Font.Underline = -4142
$wsZhCn.Range("I4").Font.Color = 15570276
$wsZhCn.Range("J4").Value = "307fbd6b-89cd-4f2d-882b-87102e101d0d.af39814767e8b6f387b253361b2e6801f0f3b45c.zh-cn.xlf"
$wsZhCn.Range("K4").Value = "2016-08-12 12:32:34"

# Row 5 - a69741f6-eb82-413a-a776-326854e80aea
$wsZhCn.Range("C5").Value = $HANDED_BACK
$wsZhCn.Range("I5").Value = "a69741f6-eb82-413a-a776-326854e80aea.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7c748e17e7c44fbaa2203baa9dbeb7529a5ddb67/e2e/a69741f6-eb82-413a-a776-326854e80aea.md") | Out-Null
$wsZhCn.Range("I5").Font.Underline = -4142
$wsZhCn.Range("I5").Font.Color = 15570276
$wsZhCn.Range("J5").Value = "a69741f6-eb82-413a-a776-326854e80aea.caa7f27e9c4d6de63948e0c7de56b49075f27d5e.zh-cn.xlf"
$wsZhCn.Range("K5").Value = "2016-08-12 12:32:34"

# ---------------------------------------------------------------------
# de-de sheet: update status + fill in handback file / datetime columns
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 4 - 307fbd6b-89cd-4f2d-882b-87102e101d0d
$wsDeDe.Range("C4").Value = $HANDED_BACK
$wsDeDe.Range("I4").Value = "307fbd6b-89cd-4f2d-882b-87102e101d0d.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1b37a371caff86262dfd1141f5fe4972d8a53026/e2e/307fbd6b-89cd-4f2d-882b-87102e101d0d.md") | Out-Null
$wsDeDe.Range("I4").Font.Underline = -4142
$wsDeDe.Range("I4").Font.Color = 15570276
$wsDeDe.Range("J4").Value = "307fbd6b-89cd-4f2d-882b-87102e101d0d.af39814767e8b6f387b253361b2e6801f0f3b45c.de-de.xlf"
$wsDeDe.Range("K4").Value = "2016-08-12 12:32:43"

# Row 5 - a69741f6-eb82-413a-a776-326854e80aea
$wsDeDe.Range("C5").Value = $HANDED_BACK
$wsDeDe.Range("I5").Value = "a69741f6-eb82-413a-a776-326854e80aea.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1b37a371caff86262dfd1141f5fe4972d8a53026/e2e/a69741f6-eb82-413a-a776-326854e80aea.md") | Out-Null
$wsDeDe.Range("I5").Font.Underline = -4142
$wsDeDe.Range("I5").Font.Color = 15570276
$wsDeDe.Range("J5").Value = "a69741f6-eb82-413a-a776-326854e80aea.caa7f27e9c4d6de63948e0c7de56b49075f27d5e.de-de.xlf"
$wsDeDe.Range("K5").Value = "2016-08-12 12:32:43"

$wb.Save()
